$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing entry (C26: 16:50 -> 17:00) ---
$ws.Range("C26").Value2 = 0.70833333333333337

# --- New log entry in row 27 (clock-in only, no clock-out yet) ---
$ws.Range("B27").Value2 = 0.72916666666666663
$ws.Range("B27").NumberFormat = "h:mm"
$ws.Range("E27").Value = "start working on first order encoder!"

# --- Move the "total" row further down to make room for more entries ---
$ws.Range("A32:D32").Clear()
$ws.Range("A49").Value = "total"
$ws.Range("D49").Formula = "= SUM(D2:D30)"
$ws.Range("D49").NumberFormat = "[h]:mm"

# --- Scroll the view down a bit ---
$ws.Application.ActiveWindow.ScrollRow = 18
